# Updates the cryptos worksheet to reflect the latest scrape:
# a new "OKB" entry is inserted at row 9, shifting the rest of the
# coin list down by one row (dropping the former last row, RenderToken),
# and refreshing Price / Volume(1h) figures throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (matches the source data, which
# stores prices/percentages as text, e.g. "27.230.27") so Excel does not
# auto-coerce numeric-looking strings like "1.006" into real numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Rows 2-8: only Price (D) / Volume (E) refresh; Coin/Link unchanged.
$ws.Range("D2").Value = "27.232.95"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").Value = "1.820.42"
$ws.Range("E3").Value = "  -2.13%  "

Set-TextValue $ws.Range("D4") "1.006"
$ws.Range("E4").Value = "  -1.41%  "

Set-TextValue $ws.Range("D5") "314.27"
$ws.Range("E5").Value = "  -2.06%  "

Set-TextValue $ws.Range("D6") "1.005"
$ws.Range("E6").Value = "  -1.31%  "

Set-TextValue $ws.Range("D7") "0.4268"
$ws.Range("E7").Value = "  -2.36%  "

$ws.Range("E8").Value = "  -2.86%  "

# New coin inserted at row 9 (OKB); rows 9-51 hold the full post-shift
# content (Coin, Link, Price, Volume(1h)).
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D9") "45.86"
$ws.Range("E9").Value = "  -2.52%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D10") "0.07223"
$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D11") "0.8603"
$ws.Range("E11").Value = "  -2.72%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D12") "20.95"
$ws.Range("E12").Value = "  -2.83%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.827.35"
$ws.Range("E13").Value = "  -1.78%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "6.663"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D15") "0.07113"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D16") "5.299"
$ws.Range("E16").Value = "  -3.57%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D17") "88.15"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D18") "1.007"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D19") "0.000008857"
$ws.Range("E19").Value = "  -1.93%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D20") "1.004"
$ws.Range("E20").Value = "  -1.40%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D21") "15.02"
$ws.Range("E21").Value = "  -2.92%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.258.26"
$ws.Range("E22").Value = "  -1.67%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D23") "5.141"
$ws.Range("E23").Value = "  -2.73%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D24") "10.88"
$ws.Range("E24").Value = "  -2.49%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.047.31"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D26") "2.002"
$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D27") "153.17"
$ws.Range("E27").Value = "  -2.53%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "18.34"
$ws.Range("E28").Value = "  -1.92%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "2.114"
$ws.Range("E29").Value = "  +6.25%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D30") "5.214"
$ws.Range("E30").Value = "  -3.90%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D31") "116.14"
$ws.Range("E31").Value = "  -4.30%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D32") "0.08887"
$ws.Range("E32").Value = "  -1.77%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D33") "1.193"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D34") "0.7566"
$ws.Range("E34").Value = "  -1.76%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D35") "4.439"
$ws.Range("E35").Value = "  -2.67%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D36") "2.821"
$ws.Range("E36").Value = "  -7.07%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D37") "1.004"
$ws.Range("E37").Value = "  -1.48%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.112"
$ws.Range("E38").Value = "  -2.48%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.01962"
$ws.Range("E39").Value = "  -0.82%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D40") "0.05265"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D41") "2.905"
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D42") "7.124"
$ws.Range("E42").Value = "  +2.37%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1675"
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D44") "0.5021"
$ws.Range("E44").Value = "  -3.08%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D45") "8.597"
$ws.Range("E45").Value = "  -1.30%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "10.59"
$ws.Range("E46").Value = "  -1.75%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D47") "106.68"
$ws.Range("E47").Value = "  -3.23%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D48") "0.4698"
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D49") "1.005"
$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.06386"
$ws.Range("E50").Value = "  -1.40%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.660"
$ws.Range("E51").Value = "  -2.99%  "
